$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.109.23"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.60"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.54"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6282"
$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07492"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2923"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("E10").Value = "  +3.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07686"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.835.53"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.010"
$ws.Range("E13").Value = "  +1.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6666"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.76"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009373"
$ws.Range("E16").Value = "  -8.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.982"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.112.97"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.081.19"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.59"
$ws.Range("E20").Value = "  +2.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "223.03"
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("E23").Value = "  -1.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.96"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1394"
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.502"
$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.88"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05642"
$ws.Range("E30").Value = "  +8.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.154"
$ws.Range("E31").Value = "  +1.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.083"
$ws.Range("E32").Value = "  +1.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.210"
$ws.Range("E33").Value = "  +1.96%  "

$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7411"
$ws.Range("E35").Value = "  +0.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.672"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.760"
$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.220.72"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01779"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.538"
$ws.Range("E41").Value = "  +3.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8917"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.79"
$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.979.80"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.80"
$ws.Range("E46").Value = "  +2.88%  "

$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5093"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4077"
$ws.Range("E49").Value = "  +1.09%  "

$ws.Range("E50").Value = "  +7.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.010"
$ws.Range("E51").Value = "  +1.83%  "
